$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'316.42"
$ws.Range("E2").Value = "'1.45%"
$ws.Range("G2").Value = "'23"
$ws.Range("D3").Value = "'37.99"
$ws.Range("E3").Value = "'1.51%"
$ws.Range("G3").Value = "'23"
$ws.Range("D4").Value = "'5.194"
$ws.Range("E4").Value = "'1.57%"
$ws.Range("G4").Value = "'23"
$ws.Range("D5").Value = "'0.08001"
$ws.Range("E5").Value = "'1.45%"
$ws.Range("G5").Value = "'23"
$ws.Range("D6").Value = "'4.483"
$ws.Range("E6").Value = "'1.38%"
$ws.Range("G6").Value = "'23"
$ws.Range("D7").Value = "'8.503"
$ws.Range("E7").Value = "'2.94%"
$ws.Range("G7").Value = "'23"
$ws.Range("D8").Value = "'1.931"
$ws.Range("E8").Value = "'1.27%"
$ws.Range("G8").Value = "'23"
$ws.Range("D9").Value = "'2.956"
$ws.Range("E9").Value = "'4.68%"
$ws.Range("G9").Value = "'23"
$ws.Range("D10").Value = "'0.9427"
$ws.Range("E10").Value = "'2.45%"
$ws.Range("G10").Value = "'23"
$ws.Range("D11").Value = "'0.1298"
$ws.Range("E11").Value = "'8.31%"
$ws.Range("G11").Value = "'23"
$ws.Range("D12").Value = "'0.1940"
$ws.Range("E12").Value = "'0.75%"
$ws.Range("G12").Value = "'23"
$ws.Range("D13").Value = "'0.09074"
$ws.Range("E13").Value = "'0.24%"
$ws.Range("G13").Value = "'23"
$ws.Range("D14").Value = "'0.03351"
$ws.Range("E14").Value = "'1.08%"
$ws.Range("G14").Value = "'23"
$ws.Range("D15").Value = "'0.09536"
$ws.Range("E15").Value = "'-0.62%"
$ws.Range("G15").Value = "'23"
$ws.Range("D16").Value = "'0.001398"
$ws.Range("E16").Value = "'1.35%"
$ws.Range("G16").Value = "'23"
$ws.Range("D17").Value = "'0.006437"
$ws.Range("E17").Value = "'7.77%"
$ws.Range("G17").Value = "'23"
$ws.Range("D18").Value = "'3.396"
$ws.Range("E18").Value = "'-4.43%"
$ws.Range("G18").Value = "'23"
$ws.Range("E19").Value = "'2.12%"
$ws.Range("G19").Value = "'23"
$ws.Range("D20").Value = "'6.576"
$ws.Range("E20").Value = "'25.88%"
$ws.Range("G20").Value = "'23"
$ws.Range("D21").Value = "'0.1314"
$ws.Range("E21").Value = "'2.32%"
$ws.Range("G21").Value = "'23"
$ws.Range("D22").Value = "'0.2424"
$ws.Range("E22").Value = "'-6.32%"
$ws.Range("G22").Value = "'23"
$ws.Range("D23").Value = "'0.04371"
$ws.Range("E23").Value = "'0.18%"
$ws.Range("G23").Value = "'23"
$ws.Range("D24").Value = "'0.001229"
$ws.Range("E24").Value = "'-1.40%"
$ws.Range("G24").Value = "'23"
$ws.Range("D25").Value = "'0.004260"
$ws.Range("E25").Value = "'-8.72%"
$ws.Range("G25").Value = "'23"
$ws.Range("E26").Value = "'-2.06%"
$ws.Range("G26").Value = "'23"
$ws.Range("D27").Value = "'0.0003988"
$ws.Range("E27").Value = "'0.13%"
$ws.Range("G27").Value = "'23"
$ws.Range("G28").Value = "'23"
$ws.Range("G29").Value = "'23"
$ws.Range("G30").Value = "'23"
$ws.Range("G31").Value = "'23"
$ws.Range("G32").Value = "'23"
$ws.Range("G33").Value = "'23"
$ws.Range("G34").Value = "'23"
$ws.Range("G35").Value = "'23"
$ws.Range("G36").Value = "'23"
$ws.Range("G37").Value = "'23"
$ws.Range("G38").Value = "'23"
$ws.Range("D39").Value = "'0.02380"
$ws.Range("E39").Value = "'5.40%"
$ws.Range("G39").Value = "'23"
$ws.Range("D40").Value = "'0.05163"
$ws.Range("E40").Value = "'1.69%"
$ws.Range("G40").Value = "'23"
$ws.Range("D41").Value = "'0.007689"
$ws.Range("E41").Value = "'3.35%"
$ws.Range("G41").Value = "'23"
$ws.Range("E42").Value = "'3.49%"
$ws.Range("G42").Value = "'23"
$ws.Range("D43").Value = "'0.008655"
$ws.Range("E43").Value = "'-4.23%"
$ws.Range("G43").Value = "'23"
$ws.Range("D44").Value = "'0.002109"
$ws.Range("E44").Value = "'8.36%"
$ws.Range("G44").Value = "'23"
$ws.Range("D45").Value = "'0.008026"
$ws.Range("E45").Value = "'-13.24%"
$ws.Range("G45").Value = "'23"
$ws.Range("D46").Value = "'0.00006474"
$ws.Range("E46").Value = "'-1.37%"
$ws.Range("G46").Value = "'23"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("G47").Value = "'23"
$ws.Range("D48").Value = "'0.002865"
$ws.Range("E48").Value = "'-13.47%"
$ws.Range("G48").Value = "'23"
$ws.Range("D49").Value = "'0.001688"
$ws.Range("E49").Value = "'69.17%"
$ws.Range("G49").Value = "'23"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("G50").Value = "'23"
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("G51").Value = "'23"
